$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the first sheet
$ws.Name = "Admin_SZ_Portuguese"

# Add new row 4 with the "Rooms" test case data
$ws.Range("A4").Value = "fn_verifyAdditionOfRoom"
$ws.Range("C4").Value = "RoomsList_Title"
$ws.Range("E4").Value = "Message_Text1"
$ws.Range("H4").Value = " salvo com sucesso"
$ws.Range("D4").Value = "Lista de quartos"
$ws.Range("G4").Value = "Message_Text2"
$ws.Range("F4").Value = "'- Quarto "

# Fix up the existing row 3, column F value: change "-  Amenidade/serviço" (double space)
# to "- Amenidade/serviço " (single space, trailing space). Leading apostrophe forces the
# "quote prefix" (text starting with a special char) formatting, matching style s="2".
$ws.Range("F3").Value = "'- Amenidade/serviço "

# Update the active selection to F3
$ws.Range("F3").Select()
